$d = $word.ActiveDocument

# The fee-table ("Rechnungspositionen") is the 3rd table in the document;
# its last row ("Total") holds the literal "0.00" placeholder that needs
# to become the {{GEBUEHREN_TOTAL}} merge field.
$t = $d.Tables.Item(3)
$cell = $t.Cell($t.Rows.Count, $t.Columns.Count)
$cell.Range.Text = "{{GEBUEHREN_TOTAL}}"
